# p2_noControls.xlsx regression-output table got refreshed with slightly
# different coefficient estimates (a rerun of the underlying regression).
# Only the standard-error / p-value cells that actually moved need updating;
# everything else (headers, A-column coefficients, C2) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 265.39269635801594
$ws.Range("B3").Value = 330.4695657555593
$ws.Range("C3").Value = 0.00045073939597006785
